# Auto-generated script to update Sheet1 rows 2-65 (A:D)
# with the new municipio/assentamento/codigo data, and extend the
# sheet's dimension from A1:D28 to A1:D65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('BANDEIRANTES', 'BANDEIRANTES', 'PR0260000', 'Sim')
    ,@('BITURUNA', '12 DE ABRIL', 'PR0178000', 'Sim')
    ,@('BITURUNA', 'CRICIUMINHA', 'PR0191000', 'Sim')
    ,@('BITURUNA', 'ETIENE', 'PR0062000', 'Sim')
    ,@('BITURUNA', 'RONDON III', 'PR0063000', 'Sim')
    ,@('BITURUNA', 'SANTA BÁRBARA', 'PR0175000', 'Sim')
    ,@('BITURUNA', 'SONHO DE ROSE', 'PR0383000', 'Sim')
    ,@('CANDOI', '8 DE OUTUBRO', 'PR0328000', 'Sim')
    ,@('CANDOI', 'COLÔNIA SÃO JOÃO BATISTA', 'PR0015000', 'Sim')
    ,@('CANDOI', 'MATA DO CAVERNOSO', 'PR0258000', 'Sim')
    ,@('CANDOI', 'UNIÃO SÃO PEDRO', 'PR0385000', 'Sim')
    ,@('CANDOI', 'ÁGUAS DE SANTA CLARA', 'PR0237000', 'Sim')
    ,@('GENERAL CARNEIRO', 'COLINA VERDE', 'PR0141000', 'Sim')
    ,@('GENERAL CARNEIRO', 'RECANTO BONITO', 'PR0240000', 'Sim')
    ,@('GENERAL CARNEIRO', 'RIO CATEQUESE', 'PR0255000', 'Sim')
    ,@('GUARAPUAVA', '13 DE NOVEMBRO', 'PR0418000', 'Sim')
    ,@('GUARAPUAVA', 'EUROPA', 'PR0234000', 'Sim')
    ,@('GUARAPUAVA', 'FAZENDA BANANAS', 'PR0246000', 'Sim')
    ,@('GUARAPUAVA', 'FAZENDA CAROLINA', 'PR0154000', 'Sim')
    ,@('GUARAPUAVA', 'NOVA GERAÇÃO', 'PR0325000', 'Sim')
    ,@('GUARAPUAVA', 'SÃO PEDRO', 'PR0077000', 'Sim')
    ,@('INACIO MARTINS', 'BOM RETIRO', 'PR0074000', 'Sim')
    ,@('INACIO MARTINS', 'EVANDRO FRANCISCO', 'PR0375000', 'Sim')
    ,@('INACIO MARTINS', 'FAXINAL DOS RODRIGUES', 'PR0117000', 'Sim')
    ,@('INACIO MARTINS', 'JOSÉ DIAS', 'PR0122000', 'Sim')
    ,@('LAPA', 'CONTESTADO', 'PR0259000', 'Sim')
    ,@('LARANJAL', 'CONQUISTA CAMPONESA', 'PR0311000', 'Sim')
    ,@('LARANJAL', 'FAZENDA CHAPADÃO', 'PR0110000', 'Sim')
    ,@('LARANJAL', 'PEDRA BRANCA', 'PR0142000', 'Sim')
    ,@('MANGUEIRINHA', 'ANJO DA GUARDA I', 'PR0330000', 'Sim')
    ,@('MANGUEIRINHA', 'ESPERANÇA VIVA', 'PR0329000', 'Sim')
    ,@('MANGUEIRINHA', 'FAZENDA COVOZINHO', 'PR0018000', 'Sim')
    ,@('MANGUEIRINHA', 'FAZENDA MACHADO', 'PR0059000', 'Sim')
    ,@('MANGUEIRINHA', 'NATAL DA ESPERANÇA', 'PR0248000', 'Sim')
    ,@('MANGUEIRINHA', 'SÃO JOÃO MARIA', 'PR0308000', 'Sim')
    ,@('MANGUEIRINHA', 'TRÊS CAPÕES', 'PR0016000', 'Sim')
    ,@('MANGUEIRINHA', 'VITÓRIA', 'PR0285000', 'Sim')
    ,@('MANGUEIRINHA', 'VITÓRIA DA UNIÃO DO PARANÁ', 'PR0017000', 'Sim')
    ,@('MANOEL RIBAS', 'NOVA ITAÚNA', 'PR0261000', 'Sim')
    ,@('MARMELEIRO', 'ANJO DA GUARDA', 'PR0145000', 'Sim')
    ,@('MARMELEIRO', 'EDUARDO RADUAN', 'PR0207000', 'Sim')
    ,@('MARMELEIRO', 'FAZENDA PERSEVERANÇA', 'PR0028000', 'Sim')
    ,@('MARMELEIRO', 'KARL MARX', 'PR0333000', 'Sim')
    ,@('MARMELEIRO', 'SÃO JORGE', 'PR0133000', 'Sim')
    ,@('PINHAO', 'FAXINAL DOS RIBEIROS QUINHÃO 1- A', 'PR0048000', 'Sim')
    ,@('PINHAO', 'FAXINAL DOS SILVÉRIOS', 'PR0073000', 'Sim')
    ,@('RESERVA', 'CAIÇARA', 'PR0106000', 'Sim')
    ,@('RESERVA', 'CRICIUMA', 'PR0058000', 'Sim')
    ,@('RESERVA', 'FAZENDA BARREIROS', 'PR0263000', 'Sim')
    ,@('RESERVA', 'FIO DE OURO', 'PR0032000', 'Sim')
    ,@('RESERVA', 'PAINEIRA', 'PR0302000', 'Sim')
    ,@('RESERVA', 'RECANTO DA AMIZADE', 'PR0249000', 'Sim')
    ,@('RESERVA', 'RENASCENÇA', 'PR0204000', 'Sim')
    ,@('RESERVA', 'RODEIO', 'PR0236000', 'Sim')
    ,@('RESERVA', 'SANTA HELENA', 'PR0031000', 'Sim')
    ,@('RESERVA', 'SINHÁ ANA', 'PR0107000', 'Sim')
    ,@('TIBAGI', 'BOA VISTA', 'PR0149000', 'Sim')
    ,@('TIBAGI', 'DONA TONIA', 'PR0244000', 'Sim')
    ,@('TIBAGI', 'FAZENDA RETIRO', 'PR0079000', 'Sim')
    ,@('TIBAGI', 'MENINO JESUS', 'PR0161000', 'Sim')
    ,@('TIBAGI', 'NOSSA SENHORA APARECIDA', 'PR0407000', 'Sim')
    ,@('TIBAGI', 'RANCHO ALEGRE', 'PR0228000', 'Sim')
    ,@('TIBAGI', 'RINCÃO', 'PR0200000', 'Sim')
    ,@('TIBAGI', 'VASTO HORIZONTE', 'PR0026000', 'Sim')
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
